$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_12_2_1"
$ws.Range("B2").Value = 0.08784013135839541
$ws.Range("C2").Value = 0.1960313320589125
$ws.Range("D2").Value = -7.039412453450664
$ws.Range("E2").Value = -0.0989564180031457
$ws.Range("F2").Value = 1.009492158889771
$ws.Range("G2").Value = 1.384710073471069
$ws.Range("H2").Value = 0.9096595048904419
$ws.Range("I2").Value = 1.161158084869385

$ws.Range("A3").Value = "model_12_2_3"
$ws.Range("B3").Value = 0.1375732173525326
$ws.Range("C3").Value = -0.1105779537300349
$ws.Range("D3").Value = -7.659803193464159
$ws.Range("E3").Value = -0.3948228453951272
$ws.Range("F3").Value = 0.954452395439148
$ws.Range("G3").Value = 1.912796258926392
$ws.Range("H3").Value = 0.9798567891120911
$ws.Range("I3").Value = 1.47377073764801

$ws.Range("A4").Value = "model_12_2_2"
$ws.Range("B4").Value = 0.1634858989596194
$ws.Range("C4").Value = 0.2277415947864321
$ws.Range("D4").Value = -7.004864560493413
$ws.Range("E4").Value = -0.06984884197252339
$ws.Range("F4").Value = 0.9257745742797852
$ws.Range("G4").Value = 1.33009397983551
$ws.Range("H4").Value = 0.9057504534721375
$ws.Range("I4").Value = 1.13040292263031

$ws.Range("A5").Value = "model_12_2_4"
$ws.Range("B5").Value = 0.168960404665644
$ws.Range("C5").Value = -0.007253844271601029
$ws.Range("D5").Value = -7.825180395176012
$ws.Range("E5").Value = -0.3139827439344567
$ws.Range("F5").Value = 0.919715940952301
$ws.Range("G5").Value = 1.73483669757843
$ws.Range("H5").Value = 0.9985692501068115
$ws.Range("I5").Value = 1.388355016708374

$ws.Range("A6").Value = "model_12_2_0"
$ws.Range("B6").Value = 0.2490739897418898
$ws.Range("C6").Value = 0.8003617340255746
$ws.Range("D6").Value = -0.7064691575528568
$ws.Range("E6").Value = 0.7417150706239545
$ws.Range("F6").Value = 0.8310537934303284
$ws.Range("G6").Value = 0.3438456058502197
$ws.Range("H6").Value = 0.1930869966745377
$ws.Range("I6").Value = 0.2729040384292603

$ws.Range("A7").Value = "model_12_2_5"
$ws.Range("B7").Value = 0.2597434543301836
$ws.Range("C7").Value = 0.226854744025784
$ws.Range("D7").Value = -9.032607178025197
$ws.Range("E7").Value = -0.1728013225634195
$ws.Range("F7").Value = 0.8192458748817444
$ws.Range("G7").Value = 1.331621527671814
$ws.Range("H7").Value = 1.135189533233643
$ws.Range("I7").Value = 1.239182591438293

$ws.Range("A8").Value = "model_12_2_6"
$ws.Range("B8").Value = 0.2736879263779448
$ws.Range("C8").Value = 0.2128805143831993
$ws.Range("D8").Value = -8.966855642516753
$ws.Range("E8").Value = -0.1815474111484305
$ws.Range("F8").Value = 0.8038134574890137
$ws.Range("G8").Value = 1.355689883232117
$ws.Range("H8").Value = 1.127749800682068
$ws.Range("I8").Value = 1.24842357635498

$ws.Range("A9").Value = "model_12_2_7"
$ws.Range("B9").Value = 0.3204415081996024
$ws.Range("C9").Value = 0.3522019101888761
$ws.Range("D9").Value = -8.853129004141731
$ws.Range("E9").Value = -0.05558416119851861
$ws.Range("F9").Value = 0.752070963382721
$ws.Range("G9").Value = 1.1157306432724
$ws.Range("H9").Value = 1.11488151550293
$ws.Range("I9").Value = 1.115330815315247

$ws.Range("A10").Value = "model_12_2_8"
$ws.Range("B10").Value = 0.3245167877805639
$ws.Range("C10").Value = 0.3719377945417058
$ws.Range("D10").Value = -8.967850795419603
$ws.Range("E10").Value = -0.04433371474718317
$ws.Range("F10").Value = 0.7475609183311462
$ws.Range("G10").Value = 1.081738710403442
$ws.Range("H10").Value = 1.127862453460693
$ws.Range("I10").Value = 1.103443741798401

$ws.Range("A11").Value = "model_12_2_10"
$ws.Range("B11").Value = 0.3431264298986386
$ws.Range("C11").Value = 0.4162059970136496
$ws.Range("D11").Value = -9.044877514968526
$ws.Range("E11").Value = -0.01001332423428702
$ws.Range("F11").Value = 0.7269654870033264
$ws.Range("G11").Value = 1.005493640899658
$ws.Range("H11").Value = 1.136577963829041
$ws.Range("I11").Value = 1.067180752754211

$ws.Range("A12").Value = "model_12_2_9"
$ws.Range("B12").Value = 0.3443319929828985
$ws.Range("C12").Value = 0.4260090626480175
$ws.Range("D12").Value = -9.045639577965616
$ws.Range("E12").Value = -0.001590934982644399
$ws.Range("F12").Value = 0.7256312966346741
$ws.Range("G12").Value = 0.9886094331741333
$ws.Range("H12").Value = 1.136664152145386
$ws.Range("I12").Value = 1.058281660079956

$ws.Range("A13").Value = "model_12_2_11"
$ws.Range("B13").Value = 0.3466515746479421
$ws.Range("C13").Value = 0.4193325080764867
$ws.Range("D13").Value = -8.869577016099955
$ws.Range("E13").Value = 0.001519530737263297
$ws.Range("F13").Value = 0.723064124584198
$ws.Range("G13").Value = 1.00010871887207
$ws.Range("H13").Value = 1.116742610931396
$ws.Range("I13").Value = 1.054995179176331

$ws.Range("A14").Value = "model_12_2_12"
$ws.Range("B14").Value = 0.34936318046822
$ws.Range("C14").Value = 0.4231329618989638
$ws.Range("D14").Value = -8.761985747852822
$ws.Range("E14").Value = 0.01022118026580965
$ws.Range("F14").Value = 0.7200632095336914
$ws.Range("G14").Value = 0.9935629963874817
$ws.Range("H14").Value = 1.104568719863892
$ws.Range("I14").Value = 1.045800924301147

$ws.Range("A15").Value = "model_12_2_13"
$ws.Range("B15").Value = 0.3530538810332304
$ws.Range("C15").Value = 0.4290445489790222
$ws.Range("D15").Value = -8.612831913062367
$ws.Range("E15").Value = 0.02283909458266009
$ws.Range("F15").Value = 0.7159786820411682
$ws.Range("G15").Value = 0.9833812713623047
$ws.Range("H15").Value = 1.087692022323608
$ws.Range("I15").Value = 1.032468914985657

$ws.Range("A16").Value = "model_12_2_14"
$ws.Range("B16").Value = 0.3533686008114044
$ws.Range("C16").Value = 0.431742173039097
$ws.Range("D16").Value = -8.638028200150883
$ws.Range("E16").Value = 0.02389774744278839
$ws.Range("F16").Value = 0.7156304121017456
$ws.Range("G16").Value = 0.9787349700927734
$ws.Range("H16").Value = 1.090543031692505
$ws.Range("I16").Value = 1.031350255012512

$ws.Range("A17").Value = "model_12_2_15"
$ws.Range("B17").Value = 0.353672031788087
$ws.Range("C17").Value = 0.4335885947651786
$ws.Range("D17").Value = -8.649026073955836
$ws.Range("E17").Value = 0.0249365123769496
$ws.Range("F17").Value = 0.7152945399284363
$ws.Range("G17").Value = 0.9755547642707825
$ws.Range("H17").Value = 1.091787338256836
$ws.Range("I17").Value = 1.030252814292908

$ws.Range("A18").Value = "model_12_2_16"
$ws.Range("B18").Value = 0.353672031788087
$ws.Range("C18").Value = 0.4335885947651786
$ws.Range("D18").Value = -8.649026073955836
$ws.Range("E18").Value = 0.0249365123769496
$ws.Range("F18").Value = 0.7152945399284363
$ws.Range("G18").Value = 0.9755547642707825
$ws.Range("H18").Value = 1.091787338256836
$ws.Range("I18").Value = 1.030252814292908

$ws.Range("A19").Value = "model_12_2_17"
$ws.Range("B19").Value = 0.353672031788087
$ws.Range("C19").Value = 0.4335885947651786
$ws.Range("D19").Value = -8.649026073955836
$ws.Range("E19").Value = 0.0249365123769496
$ws.Range("F19").Value = 0.7152945399284363
$ws.Range("G19").Value = 0.9755547642707825
$ws.Range("H19").Value = 1.091787338256836
$ws.Range("I19").Value = 1.030252814292908

$ws.Range("A20").Value = "model_12_2_18"
$ws.Range("B20").Value = 0.353672031788087
$ws.Range("C20").Value = 0.4335885947651786
$ws.Range("D20").Value = -8.649026073955836
$ws.Range("E20").Value = 0.0249365123769496
$ws.Range("F20").Value = 0.7152945399284363
$ws.Range("G20").Value = 0.9755547642707825
$ws.Range("H20").Value = 1.091787338256836
$ws.Range("I20").Value = 1.030252814292908

$ws.Range("A21").Value = "model_12_2_19"
$ws.Range("B21").Value = 0.353672031788087
$ws.Range("C21").Value = 0.4335885947651786
$ws.Range("D21").Value = -8.649026073955836
$ws.Range("E21").Value = 0.0249365123769496
$ws.Range("F21").Value = 0.7152945399284363
$ws.Range("G21").Value = 0.9755547642707825
$ws.Range("H21").Value = 1.091787338256836
$ws.Range("I21").Value = 1.030252814292908

$ws.Range("A22").Value = "model_12_2_20"
$ws.Range("B22").Value = 0.353672031788087
$ws.Range("C22").Value = 0.4335885947651786
$ws.Range("D22").Value = -8.649026073955836
$ws.Range("E22").Value = 0.0249365123769496
$ws.Range("F22").Value = 0.7152945399284363
$ws.Range("G22").Value = 0.9755547642707825
$ws.Range("H22").Value = 1.091787338256836
$ws.Range("I22").Value = 1.030252814292908

$ws.Range("A23").Value = "model_12_2_21"
$ws.Range("B23").Value = 0.353672031788087
$ws.Range("C23").Value = 0.4335885947651786
$ws.Range("D23").Value = -8.649026073955836
$ws.Range("E23").Value = 0.0249365123769496
$ws.Range("F23").Value = 0.7152945399284363
$ws.Range("G23").Value = 0.9755547642707825
$ws.Range("H23").Value = 1.091787338256836
$ws.Range("I23").Value = 1.030252814292908

$ws.Range("A24").Value = "model_12_2_22"
$ws.Range("B24").Value = 0.353672031788087
$ws.Range("C24").Value = 0.4335885947651786
$ws.Range("D24").Value = -8.649026073955836
$ws.Range("E24").Value = 0.0249365123769496
$ws.Range("F24").Value = 0.7152945399284363
$ws.Range("G24").Value = 0.9755547642707825
$ws.Range("H24").Value = 1.091787338256836
$ws.Range("I24").Value = 1.030252814292908

$ws.Range("A25").Value = "model_12_2_23"
$ws.Range("B25").Value = 0.353672031788087
$ws.Range("C25").Value = 0.4335885947651786
$ws.Range("D25").Value = -8.649026073955836
$ws.Range("E25").Value = 0.0249365123769496
$ws.Range("F25").Value = 0.7152945399284363
$ws.Range("G25").Value = 0.9755547642707825
$ws.Range("H25").Value = 1.091787338256836
$ws.Range("I25").Value = 1.030252814292908

$ws.Range("A26").Value = "model_12_2_24"
$ws.Range("B26").Value = 0.353672031788087
$ws.Range("C26").Value = 0.4335885947651786
$ws.Range("D26").Value = -8.649026073955836
$ws.Range("E26").Value = 0.0249365123769496
$ws.Range("F26").Value = 0.7152945399284363
$ws.Range("G26").Value = 0.9755547642707825
$ws.Range("H26").Value = 1.091787338256836
$ws.Range("I26").Value = 1.030252814292908
